$wb = $excel.ActiveWorkbook

# --- Text change: "Ready for handoff" -> "In Translation" ---
# This shared string is used on the Overview sheet (zh-cn / de-de status
# columns) as well as on each per-locale sheet's "Status" column.
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = "In Translation"
$wsOverview.Range("F2").Value = "In Translation"

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = "In Translation"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = "In Translation"

# --- Column width change ---
# The status text got shorter ("Ready for handoff" -> "In Translation"),
# so the status columns narrow accordingly (was 17.2159881591797, now
# 13.4101845877511 in the source workbook).
$wsOverview.Range("E:E").ColumnWidth = 12.5
$wsOverview.Range("F:F").ColumnWidth = 12.5
$wsZhCn.Range("C:C").ColumnWidth = 12.5
$wsDeDe.Range("C:C").ColumnWidth = 12.5
